$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $prefixedText) {
    $ws.Range("A1").Value = $prefixedText
    $ws.Range("A1").Copy()
    $ws.Range($cell).PasteSpecial()
    $ws.Range("A1").Clear()
}

$ws.Range('D2').Value = '27.682.53'
$ws.Range('E2').Value = '  -0.51%  '

$ws.Range('D3').Value = '1.584.82'
$ws.Range('E3').Value = '  -2.48%  '

$ws.Range('E4').Value = '  +0.72%  '

Set-TextValue 'D5' "'206.82"
$ws.Range('E5').Value = '  -1.85%  '

$ws.Range('E6').Value = '  -3.22%  '

$ws.Range('E7').Value = '  +0.76%  '

Set-TextValue 'D8' "'22.15"
$ws.Range('E8').Value = '  -4.58%  '

$ws.Range('E9').Value = '  -1.50%  '

Set-TextValue 'D10' "'0.0591"
$ws.Range('E10').Value = '  -2.90%  '

Set-TextValue 'D11' "'0.0866"
$ws.Range('E11').Value = '  -1.39%  '

$ws.Range('D12').Value = '1.809.86'
$ws.Range('E12').Value = '  -2.52%  '

$ws.Range('D13').Value = '1.621.82'
$ws.Range('E13').Value = '  -0.10%  '

Set-TextValue 'D14' "'3.84"
$ws.Range('E14').Value = '  -4.12%  '

Set-TextValue 'D15' "'0.529"
$ws.Range('E15').Value = '  -5.09%  '

Set-TextValue 'D16' "'63.49"
$ws.Range('E16').Value = '  -2.30%  '

$ws.Range('D17').Value = '27.629.68'
$ws.Range('E17').Value = '  -0.79%  '

Set-TextValue 'D18' "'219.54"
$ws.Range('E18').Value = '  -3.60%  '

$ws.Range('D19').Value = '0.0₃0693'
$ws.Range('E19').Value = '  -3.49%  '

Set-TextValue 'D20' "'7.32"
$ws.Range('E20').Value = '  -3.32%  '

$ws.Range('E21').Value = '  +0.75%  '

Set-TextValue 'D22' "'4.13"
$ws.Range('E22').Value = '  -4.33%  '

Set-TextValue 'D23' "'9.54"
$ws.Range('E23').Value = '  -3.94%  '

$ws.Range('E24').Value = '  -3.51%  '

Set-TextValue 'D25' "'153.83"
$ws.Range('E25').Value = '  -1.02%  '

Set-TextValue 'D26' "'6.87"
$ws.Range('E26').Value = '  -0.79%  '

$ws.Range('E27').Value = '  +0.76%  '

Set-TextValue 'D28' "'15.12"
$ws.Range('E28').Value = '  -2.06%  '

$ws.Range('E29').Value = '  -4.36%  '

$ws.Range('E30').Value = '  -1.97%  '

$ws.Range('E31').Value = '  -2.63%  '

Set-TextValue 'D32' "'3.21"
$ws.Range('E32').Value = '  -5.54%  '

$ws.Range('D33').Value = '1.357.91'
$ws.Range('E33').Value = '  -3.34%  '

Set-TextValue 'D34' "'2.94"
$ws.Range('E34').Value = '  -4.71%  '

$ws.Range('E35').Value = '  -4.88%  '

Set-TextValue 'D36' "'0.973"
$ws.Range('E36').Value = '  -2.85%  '

Set-TextValue 'D37' "'2.30"
$ws.Range('E37').Value = '  -0.91%  '

Set-TextValue 'D38' "'0.0168"
$ws.Range('E38').Value = '  -1.20%  '

$ws.Range('E39').Value = '  -3.34%  '

$ws.Range('E40').Value = '  -3.04%  '

$ws.Range('E41').Value = '  +0.70%  '

Set-TextValue 'D42' "'0.977"
$ws.Range('E42').Value = '  -1.98%  '

Set-TextValue 'D43' "'63.74"
$ws.Range('E43').Value = '  -2.88%  '

$ws.Range('E44').Value = '  +2.81%  '

Set-TextValue 'D45' "'1.73"
$ws.Range('E45').Value = '  -4.04%  '

Set-TextValue 'D46' "'5.18"
$ws.Range('E46').Value = '  -4.80%  '

$ws.Range('D47').Value = '1.720.62'
$ws.Range('E47').Value = '  -2.54%  '

Set-TextValue 'D48' "'88.19"
$ws.Range('E48').Value = '  +0.25%  '

$ws.Range('E49').Value = '  +12.52%  '

Set-TextValue 'D50' "'0.0968"
$ws.Range('E50').Value = '  -4.12%  '

Set-TextValue 'D51' "'0.0498"
$ws.Range('E51').Value = '  -1.04%  '

